# Weekly price-sheet update: a new weekly observation is inserted above the
# two most-recent rows (23 and 24), pushing the existing rows down by one.
# Row 23 gets filled with the new week's figures; the old row 23 becomes
# row 24 and the old row 24 becomes row 25 (unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23, shifting rows 23-24 down to 24-25.
$ws.Rows(23).Insert()

# Fill the newly inserted row 23 with this week's data.
$ws.Cells.Item(23, 1).Value  = 7
$ws.Cells.Item(23, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(23, 3).Value  = "Ñuble"
$ws.Cells.Item(23, 4).Value  = 44516
$ws.Cells.Item(23, 5).Value  = 16
$ws.Cells.Item(23, 6).Value  = 100112022
$ws.Cells.Item(23, 7).Value  = "Arveja Verde"
$ws.Cells.Item(23, 8).Value  = "Sin especificar"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 80
$ws.Cells.Item(23, 11).Value = 15000
$ws.Cells.Item(23, 12).Value = 16000
$ws.Cells.Item(23, 13).Value = 15500
$ws.Cells.Item(23, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Región del Maule"
$ws.Cells.Item(23, 16).Value = 620
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
